$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values; force text format so numeric-looking
# strings (prices with dot-grouping, tiny decimals, etc.) are preserved
# exactly as text instead of being reinterpreted as numbers/dates by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.279.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.679.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5260"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2697"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06429"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.01"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07503"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.686.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.543"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5800"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008476"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.309.38"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.919"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.07"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.715"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1236"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06561"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.358"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.578"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.564"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.658"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6185"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.398"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.380"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.104.13"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8762"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.015"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.43"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.828.03"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000114"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.79"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.146"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05270"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4304"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.041"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.39%  "
